# Apply the update described by the diff:
#  - Insert "CALCITRON 30 CAPS." as a new product row (alphabetically after
#    ANTODINE, before E-MOX)
#  - Insert "FEROGLOBIN 30 CAPS" as a new product row (alphabetically after
#    ERASTAPEX, before FLAGYL)
#  - Shift the existing rows (and the totals / footer rows) down to make
#    room, update the total, and refresh the footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: move rows 17 -> 19 and 16 -> 18 (footer + totals rows), then
# the existing product rows 15..8 down by two rows each (processing from
# the bottom up so we never overwrite a row before it has been copied).
# ---------------------------------------------------------------------

# Footer row (was row 17) -> row 19
$ws.Range("A17:Q17").Copy($ws.Range("A19:Q19"))
$ws.Rows.Item(19).RowHeight = 16.5

# Totals row (was row 16) -> row 18
$ws.Range("A16:Q16").Copy($ws.Range("A18:Q18"))
$ws.Rows.Item(18).RowHeight = 24.75

# Clear old merges that are no longer valid at rows 16/17 before we reuse
# them, and drop the old merge that copying duplicated onto rows 18/19's
# target ranges.
$ws.Range("P16:Q16").UnMerge()
$ws.Range("A17:F17").UnMerge()
$ws.Range("G17:I17").UnMerge()
$ws.Range("K17:Q17").UnMerge()

# Product rows: TICANASE (15) -> 17, سرنجات 3 سم (14) -> 16,
#               OFRAMAX (12) -> 14, LAMIFEN (11) -> 13,
#               FLAGYL (10) -> 12
# (رow 13, "سرنجات 5 سم" at old row 15, already handled above via row 17 copy target is wrong;
#  recompute precisely below in descending order.)

# old row 15 (سرنجات 5 سم) -> new row 17
$ws.Range("A15:Q15").Copy($ws.Range("A17:Q17"))
$ws.Rows.Item(17).RowHeight = 25.5

# old row 14 (سرنجات 3 سم) -> new row 16
$ws.Range("A14:Q14").Copy($ws.Range("A16:Q16"))
$ws.Rows.Item(16).RowHeight = 25.5

# old row 13 (TICANASE) -> new row 15
$ws.Range("A13:Q13").Copy($ws.Range("A15:Q15"))
$ws.Rows.Item(15).RowHeight = 24.75

# old row 12 (OFRAMAX) -> new row 14
$ws.Range("A12:Q12").Copy($ws.Range("A14:Q14"))
$ws.Rows.Item(14).RowHeight = 25.5

# old row 11 (LAMIFEN) -> new row 13
$ws.Range("A11:Q11").Copy($ws.Range("A13:Q13"))
$ws.Rows.Item(13).RowHeight = 24.75

# old row 10 (FLAGYL) -> new row 12
$ws.Range("A10:Q10").Copy($ws.Range("A12:Q12"))
$ws.Rows.Item(12).RowHeight = 25.5

# old row 9 (ERASTAPEX) -> new row 10
$ws.Range("A9:Q9").Copy($ws.Range("A10:Q10"))
$ws.Rows.Item(10).RowHeight = 24.75

# old row 8 (E-MOX) -> new row 9
$ws.Range("A8:Q8").Copy($ws.Range("A9:Q9"))
$ws.Rows.Item(9).RowHeight = 25.5

# row 7 (ANTODINE) stays in place.

# ---------------------------------------------------------------------
# Step 2: populate the two brand-new product rows (8 = CALCITRON,
# 11 = FEROGLOBIN) using row 7 as the formatting template.
# ---------------------------------------------------------------------

$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Range("A8").Value() = 2

$ws.Range("A7:Q7").Copy($ws.Range("A11:Q11"))
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Range("A11").Value() = 5

# ---------------------------------------------------------------------
# Step 3: fix up the sequence numbers in column A for every product row
# (1 .. 11) now that the table has 11 rows instead of 9.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 11; $i++) {
    $ws.Cells.Item(7 + $i, 1).Value() = $i + 1
}

# ---------------------------------------------------------------------
# Step 4: write the cell values for every product row (columns C, H, L,
# N, P, Q). Column A (seq) was already set above; B/D/E/F/G/I/J/K/M/O
# are part of merged ranges and stay blank.
# ---------------------------------------------------------------------

function Set-ProductRow($Row, $Name, $Balance, $ReorderLevel, $Price, $SalePrice, $TxnCount) {
    $ws.Range("C$Row").Value() = $Name
    $ws.Range("H$Row").Value() = $Balance
    $ws.Range("L$Row").Value() = $ReorderLevel
    $ws.Range("N$Row").Value() = $Price
    $ws.Range("P$Row").Value() = $SalePrice
    $ws.Range("Q$Row").Value() = $TxnCount
}

Set-ProductRow 7  "ANTODINE 20MG 30 F.C.TAB"          "1:0" "1" "60.00"  "19.8000" "0:1"
Set-ProductRow 8  "CALCITRON 30 CAPS."                "1:0" "1" "144.00" "47.5200" "0:1"
Set-ProductRow 9  "E-MOX 500MG 16 CAPS"                "1:0" "1" "40.00"  "20.0000" "0:1"
Set-ProductRow 10 "ERASTAPEX PLUS 40MG/12.5MG 30 TAB"  "1:0" "1" "96.00"  "96.0000" "1:0"
Set-ProductRow 11 "FEROGLOBIN 30 CAPS"                 "1:1" "1" "180.00" "90.0000" "0:1"
Set-ProductRow 12 "FLAGYL 125MG/5ML 100 ML SUSPENSION" "5:0" "1" "26.00"  "26.0000" "1:0"
Set-ProductRow 13 "LAMIFEN 250MG 14 TAB"               "0:1" "1" "112.00" "56.0000" "0:1"
Set-ProductRow 14 "OFRAMAX 1 GM I.M. VIAL"             "9:0" "1" "71.00"  "71.0000" "1:0"
Set-ProductRow 15 "TICANASE 0.05% NASAL SPRAY 12 GM"   "1:0" "1" "70.00"  "70.0000" "1:0"
Set-ProductRow 16 "سرنجات 3 سم"                         "0:0" "0" "2.00"   "6.0000"  "3:0"
Set-ProductRow 17 "سرنجات 5 سم"                         "0:0" "0" "3.00"   "3.0000"  "1:0"

# ---------------------------------------------------------------------
# Step 5: update the totals row (now row 18) and the footer row (now
# row 19, new timestamp).
# ---------------------------------------------------------------------

$ws.Range("P18").Value() = 505.32

$ws.Range("A19").Value() = "Saturday, 2 August, 2025 10:32 AM"
$ws.Range("G19").Value() = "1/1"
$ws.Range("K19").Value() = "developed by : Abdelaziz Talaat"

# ---------------------------------------------------------------------
# Step 6: (re)apply merges for the two newly created product rows and
# make sure the totals/footer merges are present at their new location.
# ---------------------------------------------------------------------

$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

$ws.Range("P18:Q18").Merge()
$ws.Range("A19:F19").Merge()
$ws.Range("G19:I19").Merge()
$ws.Range("K19:Q19").Merge()

Write-Host "Done"
